$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value to be stored as literal text (avoids Excel
# auto-converting numeric-looking strings like "0.999" or "4.20" into numbers,
# which would corrupt the display value / drop significant trailing zeros).
function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

$ws.Range("D2").Value = '67.568.07'
$ws.Range("E2").Value = '  -2.45%  '
$ws.Range("D3").Value = '2.388.10'
$ws.Range("E3").Value = '  -3.40%  '
Set-TextValue $ws.Range("D4") '0.999'
$ws.Range("E4").Value = '  -0.14%  '
Set-TextValue $ws.Range("D5") '548.78'
$ws.Range("E5").Value = '  -1.99%  '
Set-TextValue $ws.Range("D6") '156.97'
$ws.Range("E6").Value = '  -4.06%  '
$ws.Range("E7").Value = '  -0.05%  '
Set-TextValue $ws.Range("D8") '0.502'
$ws.Range("E8").Value = '  -0.46%  '
Set-TextValue $ws.Range("D9") '0.156'
$ws.Range("E9").Value = '  +2.49%  '
$ws.Range("E10").Value = '  -1.49%  '
Set-TextValue $ws.Range("D11") '0.326'
$ws.Range("E11").Value = '  -2.89%  '
Set-TextValue $ws.Range("D12") '4.72'
$ws.Range("E12").Value = '  -2.25%  '
$ws.Range("D13").Value = '67.437.33'
$ws.Range("E13").Value = '  -2.10%  '
Set-TextValue $ws.Range("D14") '0.0000167'
$ws.Range("E14").Value = '  -1.55%  '
Set-TextValue $ws.Range("D15") '22.71'
$ws.Range("E15").Value = '  -4.25%  '
Set-TextValue $ws.Range("D16") '10.27'
$ws.Range("E16").Value = '  -5.13%  '
Set-TextValue $ws.Range("D17") '328.92'
$ws.Range("E17").Value = '  -4.21%  '
Set-TextValue $ws.Range("D18") '6.73'
$ws.Range("E18").Value = '  -5.19%  '
Set-TextValue $ws.Range("D19") '3.74'
$ws.Range("E19").Value = '  -1.88%  '
$ws.Range("E20").Value = '  -0.43%  '
Set-TextValue $ws.Range("D21") '1.83'
$ws.Range("E21").Value = '  -5.72%  '
Set-TextValue $ws.Range("D22") '65.35'
$ws.Range("E22").Value = '  -2.76%  '
Set-TextValue $ws.Range("D23") '3.59'
$ws.Range("E23").Value = '  -3.08%  '
Set-TextValue $ws.Range("D24") '7.94'
$ws.Range("E24").Value = '  -3.76%  '
$ws.Range("D25").Value = '0.0₃0789'
$ws.Range("E25").Value = '  -3.75%  '
Set-TextValue $ws.Range("D26") '6.94'
$ws.Range("E26").Value = '  -3.87%  '
Set-TextValue $ws.Range("D27") '0.998'
$ws.Range("E27").Value = '  -0.20%  '
Set-TextValue $ws.Range("D28") '413.25'
$ws.Range("E28").Value = '  -6.60%  '
Set-TextValue $ws.Range("D29") '1.11'
$ws.Range("E29").Value = '  -2.97%  '
Set-TextValue $ws.Range("D30") '1.58'
$ws.Range("E30").Value = '  -2.52%  '
Set-TextValue $ws.Range("D31") '158.43'
$ws.Range("E31").Value = '  +1.54%  '
Set-TextValue $ws.Range("D32") '18.96'
$ws.Range("E32").Value = '  -0.62%  '
$ws.Range("E33").Value = '  -0.09%  '
Set-TextValue $ws.Range("D34") '17.59'
$ws.Range("E34").Value = '  -2.03%  '
Set-TextValue $ws.Range("D35") '0.104'
$ws.Range("E35").Value = '  -4.63%  '
Set-TextValue $ws.Range("D36") '0.290'
$ws.Range("E36").Value = '  -4.10%  '
Set-TextValue $ws.Range("D37") '4.20'
$ws.Range("E37").Value = '  -6.14%  '
Set-TextValue $ws.Range("D38") '1.45'
$ws.Range("E38").Value = '  -2.49%  '
Set-TextValue $ws.Range("D39") '1.05'
$ws.Range("E39").Value = '  -4.70%  '
$ws.Range("B40").Value = 'Aave'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range("D40") '127.61'
$ws.Range("E40").Value = '  -4.28%  '
$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range("D41") '3.25'
$ws.Range("E41").Value = '  -3.42%  '
Set-TextValue $ws.Range("D42") '1.92'
$ws.Range("E42").Value = '  -8.17%  '
Set-TextValue $ws.Range("D43") '0.0703'
$ws.Range("E43").Value = '  -2.35%  '
Set-TextValue $ws.Range("D44") '0.469'
$ws.Range("E44").Value = '  -3.16%  '
Set-TextValue $ws.Range("D45") '0.549'
$ws.Range("E45").Value = '  -2.80%  '
Set-TextValue $ws.Range("D46") '0.0908'
$ws.Range("E46").Value = '  -0.07%  '
$ws.Range("E47").Value = '  -1.10%  '
Set-TextValue $ws.Range("D48") '1.32'
$ws.Range("E48").Value = '  -8.49%  '
Set-TextValue $ws.Range("D49") '16.36'
$ws.Range("E49").Value = '  -3.73%  '
Set-TextValue $ws.Range("D50") '0.0425'
$ws.Range("E50").Value = '  -1.37%  '
$ws.Range("D51").Value = '0.0₆0200'
$ws.Range("E51").Value = '  -5.62%  '
